$p = $ppt.ActivePresentation

# --- Move the "Any Questions?" slide so it comes right before the
# "Let's get to business" slide (swaps their order: positions 11 and 12). ---
$p.Slides.Item(12).MoveTo(11)

# --- Fill in the body placeholder (idx=1) on the "Any Questions?" slide
# with the new GIT / website info. ---
$questionsSlide = $p.Slides.Item(11)
$nl  = [char]13
$tab = [char]9

$body = $questionsSlide.Shapes.Item("Text Placeholder 4")
$tr = $body.TextFrame.TextRange
$tr.Text = "BTW, we have GIT!" + $nl + "Check my website for downloading MATLAB and getting the GIT link" + $nl
[void]$tr.InsertAfter($tab)
$null = $tr.InsertAfter("www.nmahmoudi.ir/teaching.php")

# --- Incidental cleanups elsewhere in the deck (runs that got retyped as a
# single run with identical formatting during the same editing pass). ---

# Slide 2: "** PLAGIARISM IS ... **" was split across two runs; join them.
$slide2 = $p.Slides.Item(2)
$plagiarismShape = $slide2.Shapes.Item(1)
$plagiarismRange = $plagiarismShape.TextFrame.TextRange
$star = [string][char]0x2A
$ldq  = [string][char]0x201C
$rdq  = [string][char]0x201D
$apos = [string][char]0x2019
$plagiarismText = $star + $star + " PLAGIARISM IS " + $ldq + "NOT" + $rdq + " TOLERATED! (you" + $apos + "re better off not turning in your homework)" + $star + $star
$plagiarismFound = $plagiarismRange.Find("** PLAGIARISM ", 0)
$plagiarismChars = $plagiarismRange.Characters($plagiarismFound.Start, $plagiarismText.Length)
$plagiarismChars.Text = $plagiarismText

# Slide 4: "A de facto standard framework for academic advancements" was
# split across two runs; join them.
$slide4 = $p.Slides.Item(4)
$matlabShape = $slide4.Shapes.Item(2)
$matlabRange = $matlabShape.TextFrame.TextRange
$advancementsText = "A de facto standard framework for academic advancements"
$advancementsFound = $matlabRange.Find("A de facto standard", 0)
$advancementsChars = $matlabRange.Characters($advancementsFound.Start, $advancementsText.Length)
$advancementsChars.Text = $advancementsText
